$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header renames
$ws.Range("B1").Value = "sparsity_necessary"
$ws.Range("C1").Value = "necessary explanation rate"

# Rows 2-16 (CoDy / 1-best -> 1-delta) in column E
for ($r = 2; $r -le 16; $r++) {
    $ws.Cells.Item($r, 5).Value = "1-delta"
}

# Rows 26-36 (Greedy / recent -> temporal) in column E
for ($r = 26; $r -le 36; $r++) {
    $ws.Cells.Item($r, 5).Value = "temporal"
}
